$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric/percent-looking values are written as literal text (matching source data),
# by forcing Text number format on the Price/Volume columns before assignment.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "319.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.36%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.21%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.214"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.65%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07701"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.76%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.304"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.69%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.700"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.95%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9504"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4.57%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.425"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.22%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1268"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "9.36%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1836"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.65%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09183"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.15%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04236"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.37%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1053"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.96%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001279"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.34%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005887"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.25%"
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.004239"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.52%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.355"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.02%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3351"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.87%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.454"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "7.76%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1353"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.68%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2783"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.93%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04027"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.02%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001265"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.53%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001270"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.10%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02531"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "3.63%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05349"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.94%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007796"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.03%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1317"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.21%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007349"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.18%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001941"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.08%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007570"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.04%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3420"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "11.19%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006693"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.91%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.02%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2208"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "148.37%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004203"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.98%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.02%"
